# Update forecast values on the "Forecast Comparison" sheet (Removed Auto Arima)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# row -> @(D, E, F, G)
$values = @{
    2  = @(173, 204, 231, 273)
    3  = @(127, 151, 174, 209)
    4  = @(123, 145, 165, 197)
    5  = @(123, 146, 168, 201)
    6  = @(126, 150, 175, 213)
    7  = @(125, 149, 173, 210)
    8  = @(127, 152, 177, 217)
    9  = @(129, 155, 181, 222)
    10 = @(126, 151, 176, 215)
    11 = @(127, 152, 177, 217)
    12 = @(127, 152, 179, 221)
    13 = @(131, 158, 187, 234)
    14 = @(127, 152, 179, 222)
    15 = @(126, 152, 183, 230)
    16 = @(122, 148, 177, 221)
    17 = @(122, 147, 175, 218)
}

foreach ($row in $values.Keys) {
    $cols = $values[$row]
    $ws.Cells.Item($row, 4).Value = $cols[0]  # D
    $ws.Cells.Item($row, 5).Value = $cols[1]  # E
    $ws.Cells.Item($row, 6).Value = $cols[2]  # F
    $ws.Cells.Item($row, 7).Value = $cols[3]  # G
}
